$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force text number-format first so numeric-looking
# strings (e.g. "591.62") are stored as text, matching the original inlineStr cells,
# then restore the default "Normal" style so no stray formatting is left on the cells.
$dCells = @("D2","D3","D5","D6","D8","D12","D14","D15","D16","D17","D18","D19","D22","D23","D24","D25","D26","D30","D31","D32","D33","D35","D36","D37","D38","D41","D42","D43","D46","D48","D49","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.577.66"
$ws.Range("D3").Value = "3.629.18"
$ws.Range("D5").Value = "591.62"
$ws.Range("D6").Value = "192.12"
$ws.Range("D8").Value = "3.622.96"
$ws.Range("D12").Value = "58.24"
$ws.Range("D14").Value = "9.90"
$ws.Range("D15").Value = "4.207.27"
$ws.Range("D16").Value = "19.75"
$ws.Range("D17").Value = "3.625.78"
$ws.Range("D18").Value = "70.516.08"
$ws.Range("D19").Value = "12.68"
$ws.Range("D22").Value = "488.01"
$ws.Range("D23").Value = "19.39"
$ws.Range("D24").Value = "5.39"
$ws.Range("D25").Value = "4.47"
$ws.Range("D26").Value = "90.96"
$ws.Range("D30").Value = "33.07"
$ws.Range("D31").Value = "7.84"
$ws.Range("D32").Value = "625.63"
$ws.Range("D33").Value = "12.29"
$ws.Range("D35").Value = "66.21"
$ws.Range("D36").Value = "39.28"
$ws.Range("D37").Value = "0.413"
$ws.Range("D38").Value = "0.0₃0815"
$ws.Range("D41").Value = "3.60"
$ws.Range("D42").Value = "3.299.58"
$ws.Range("D43").Value = "3.16"
$ws.Range("D46").Value = "3.33"
$ws.Range("D48").Value = "9.18"
$ws.Range("D49").Value = "2.76"
$ws.Range("D51").Value = "1.00"

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Volume(1h) column (E) updates (plain text, percentages with padding spaces)
$ws.Range("E2").Value = "  +4.94%  "
$ws.Range("E3").Value = "  +4.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  +4.77%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +2.80%  "
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("E12").Value = "  +3.38%  "
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("E14").Value = "  +4.82%  "
$ws.Range("E15").Value = "  +5.03%  "
$ws.Range("E16").Value = "  +5.92%  "
$ws.Range("E17").Value = "  +4.96%  "
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("E19").Value = "  +4.76%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +4.16%  "
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E23").Value = "  +14.43%  "
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").Value = "  +6.63%  "
$ws.Range("E28").Value = "  +3.01%  "
$ws.Range("E29").Value = "  +5.94%  "
$ws.Range("E30").Value = "  +5.43%  "
$ws.Range("E31").Value = "  +9.80%  "
$ws.Range("E32").Value = "  +6.19%  "
$ws.Range("E33").Value = "  +4.82%  "
$ws.Range("E34").Value = "  +7.20%  "
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("E37").Value = "  +7.25%  "
$ws.Range("E38").Value = "  +6.00%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +8.82%  "
$ws.Range("E44").Value = "  +9.80%  "
$ws.Range("E45").Value = "  +5.22%  "
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("E51").Value = "  +0.01%  "
